$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row pair swaps (data was re-ordered in source feed) ---

# Row 97 (id 95): now holds the Basel vs Grasshoppers fixture (was Young Boys vs St Gallen)
$ws.Range("B97").Value2 = 6811292
$ws.Range("F97").Value2 = "Basel"
$ws.Range("G97").Value2 = "Grasshoppers"
$ws.Range("H97").Value2 = 0
$ws.Range("I97").Value2 = 1
$ws.Range("J97").Value2 = "A"
$ws.Range("K97").Value2 = 2
$ws.Range("L97").Value2 = 3.6
$ws.Range("M97").Value2 = 3.5
$ws.Range("N97").Value2 = 1.833
$ws.Range("O97").Value2 = 3.8
$ws.Range("P97").Value2 = 4.2
$ws.Range("Q97").Value2 = -0.5
$ws.Range("R97").Value2 = 1.825
$ws.Range("S97").Value2 = 2.025
$ws.Range("T97").Value2 = 2.75
$ws.Range("U97").Value2 = 1.975
$ws.Range("V97").Value2 = 1.875
$ws.Range("W97").Value2 = -1
$ws.Range("X97").Value2 = -1
$ws.Range("Y97").Value2 = 3.2
$ws.Range("Z97").Value2 = -1
$ws.Range("AA97").Value2 = 1.025
$ws.Range("AB97").Value2 = -1
$ws.Range("AC97").Value2 = 0.875

# Row 98 (id 96): now holds the Young Boys vs St Gallen fixture (was Basel vs Grasshoppers)
$ws.Range("B98").Value2 = 6811262
$ws.Range("F98").Value2 = "Young Boys"
$ws.Range("G98").Value2 = "St Gallen"
$ws.Range("H98").Value2 = 3
$ws.Range("I98").Value2 = 0
$ws.Range("J98").Value2 = "H"
$ws.Range("K98").Value2 = 2
$ws.Range("L98").Value2 = 3.8
$ws.Range("M98").Value2 = 3.3
$ws.Range("N98").Value2 = 2.3
$ws.Range("O98").Value2 = 4
$ws.Range("P98").Value2 = 2.8
$ws.Range("Q98").Value2 = -0.25
$ws.Range("R98").Value2 = 2.025
$ws.Range("S98").Value2 = 1.775
$ws.Range("T98").Value2 = 3.5
$ws.Range("U98").Value2 = 1.975
$ws.Range("V98").Value2 = 1.875
$ws.Range("W98").Value2 = 1.3
$ws.Range("X98").Value2 = -1
$ws.Range("Y98").Value2 = -1
$ws.Range("Z98").Value2 = 1.025
$ws.Range("AA98").Value2 = -1
$ws.Range("AB98").Value2 = -1
$ws.Range("AC98").Value2 = 0.875

# Row 101 (id 99): now holds the Yverdon vs Stade LausanneOuchy fixture (was FC Zurich vs Lucerne)
$ws.Range("B101").Value2 = 6810777
$ws.Range("F101").Value2 = "Yverdon Sport FC"
$ws.Range("G101").Value2 = "Stade LausanneOuchy"
$ws.Range("H101").Value2 = 2
$ws.Range("I101").Value2 = 1
$ws.Range("J101").Value2 = "H"
$ws.Range("K101").Value2 = 2.3
$ws.Range("L101").Value2 = 3.6
$ws.Range("M101").Value2 = 2.875
$ws.Range("N101").Value2 = 3
$ws.Range("O101").Value2 = 3.5
$ws.Range("P101").Value2 = 2.375
$ws.Range("Q101").Value2 = 0.25
$ws.Range("R101").Value2 = 1.8
$ws.Range("S101").Value2 = 2.05
$ws.Range("T101").Value2 = 2.5
$ws.Range("U101").Value2 = 1.925
$ws.Range("V101").Value2 = 1.925
$ws.Range("W101").Value2 = 2
$ws.Range("X101").Value2 = -1
$ws.Range("Y101").Value2 = -1
$ws.Range("Z101").Value2 = 0.8
$ws.Range("AA101").Value2 = -1
$ws.Range("AB101").Value2 = 0.925
$ws.Range("AC101").Value2 = -1

# Row 102 (id 100): now holds the FC Zurich vs Lucerne fixture (was Yverdon vs Stade LausanneOuchy)
$ws.Range("B102").Value2 = 6811260
$ws.Range("F102").Value2 = "FC Zurich"
$ws.Range("G102").Value2 = "Lucerne"
$ws.Range("H102").Value2 = 1
$ws.Range("I102").Value2 = 1
$ws.Range("J102").Value2 = "D"
$ws.Range("K102").Value2 = 1.833
$ws.Range("L102").Value2 = 3.6
$ws.Range("M102").Value2 = 4.2
$ws.Range("N102").Value2 = 1.7
$ws.Range("O102").Value2 = 3.6
$ws.Range("P102").Value2 = 4.75
$ws.Range("Q102").Value2 = -0.75
$ws.Range("R102").Value2 = 1.95
$ws.Range("S102").Value2 = 1.9
$ws.Range("T102").Value2 = 2.75
$ws.Range("U102").Value2 = 1.975
$ws.Range("V102").Value2 = 1.875
$ws.Range("W102").Value2 = -1
$ws.Range("X102").Value2 = 2.6
$ws.Range("Y102").Value2 = -1
$ws.Range("Z102").Value2 = -1
$ws.Range("AA102").Value2 = 0.8999999999999999
$ws.Range("AB102").Value2 = -1
$ws.Range("AC102").Value2 = 0.875

# Row 177 (id 175): refreshed odds for Young Boys vs Grasshoppers, id moved to 7616907
$ws.Range("B177").Value2 = 7616907
$ws.Range("E177").Value2 = 45386.64583333334
$ws.Range("F177").Value2 = "Young Boys"
$ws.Range("G177").Value2 = "Grasshoppers"
$ws.Range("K177").Value2 = 1.5
$ws.Range("L177").Value2 = 4.333
$ws.Range("M177").Value2 = 6
$ws.Range("N177").Value2 = 1.5
$ws.Range("O177").Value2 = 4.333
$ws.Range("P177").Value2 = 6.5
$ws.Range("Q177").Value2 = -1
$ws.Range("R177").Value2 = 1.8
$ws.Range("S177").Value2 = 2.05
$ws.Range("T177").Value2 = 3
$ws.Range("U177").Value2 = 2.025
$ws.Range("V177").Value2 = 1.825

# Row 178 (id 176): refreshed odds for Lausanne Sports vs St Gallen, id moved to 7616908
$ws.Range("B178").Value2 = 7616908
$ws.Range("E178").Value2 = 45386.64583333334
$ws.Range("F178").Value2 = "Lausanne Sports"
$ws.Range("G178").Value2 = "St Gallen"
$ws.Range("K178").Value2 = 2.6
$ws.Range("L178").Value2 = 3.4
$ws.Range("M178").Value2 = 2.6
$ws.Range("N178").Value2 = 2.6
$ws.Range("O178").Value2 = 3.5
$ws.Range("P178").Value2 = 2.6
$ws.Range("Q178").Value2 = 0
$ws.Range("R178").Value2 = 1.925
$ws.Range("S178").Value2 = 1.925
$ws.Range("T178").Value2 = 3
$ws.Range("U178").Value2 = 2
$ws.Range("V178").Value2 = 1.85

# Row 179 (id 177): refreshed odds for Lucerne vs Yverdon Sport FC, id moved to 7616906
$ws.Range("B179").Value2 = 7616906
$ws.Range("F179").Value2 = "Lucerne"
$ws.Range("G179").Value2 = "Yverdon Sport FC"
$ws.Range("K179").Value2 = 1.5
$ws.Range("L179").Value2 = 4.333
$ws.Range("M179").Value2 = 6
$ws.Range("N179").Value2 = 1.45
$ws.Range("O179").Value2 = 4.5
$ws.Range("P179").Value2 = 6.5
$ws.Range("Q179").Value2 = -1.25
$ws.Range("R179").Value2 = 2
$ws.Range("S179").Value2 = 1.85
$ws.Range("T179").Value2 = 3
$ws.Range("U179").Value2 = 1.95
$ws.Range("V179").Value2 = 1.9

# Rows 180 and 181 were duplicate/stale fixtures for the same matches now
# consolidated into rows 177-179 above, so remove them and shift rows up.
$ws.Rows.Item(180).Delete()
$ws.Rows.Item(180).Delete()
